$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newNote = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.02 = 15458.63 pesos`n✅ 15458.63 pesos = 3.99 = 945.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newNote

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 248.599
$ws2.Range("O10").Value = 3843
$ws2.Range("N12").Value = 3875.64
$ws2.Range("O12").Value = 237
